$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 12:20"

# Refresh per-country statistics (some countries also change row position
# in the ranking because of the new totals)
# Row 9: Iran
$ws.Cells.Item(9, 1).Value = "Iran"
$ws.Cells.Item(9, 2).Value = 38309
$ws.Cells.Item(9, 3).Value = 2901
$ws.Cells.Item(9, 4).Value = 12391
$ws.Cells.Item(9, 5).Value = 23278
$ws.Cells.Item(9, 6).Value = 3206
$ws.Cells.Item(9, 7).Value = 123
$ws.Cells.Item(9, 8).Value = 2640

# Row 10: Francia
$ws.Cells.Item(10, 1).Value = "Francia"
$ws.Cells.Item(10, 2).Value = 37575
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 5700
$ws.Cells.Item(10, 5).Value = 29561
$ws.Cells.Item(10, 6).Value = 4273
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 2314

# Row 12: Suiza
$ws.Cells.Item(12, 1).Value = "Suiza"
$ws.Cells.Item(12, 2).Value = 14352
$ws.Cells.Item(12, 3).Value = 276
$ws.Cells.Item(12, 4).Value = 1595
$ws.Cells.Item(12, 5).Value = 12475
$ws.Cells.Item(12, 6).Value = 301
$ws.Cells.Item(12, 7).Value = 18
$ws.Cells.Item(12, 8).Value = 282

# Row 16: Austria
$ws.Cells.Item(16, 1).Value = "Austria"
$ws.Cells.Item(16, 2).Value = 8411
$ws.Cells.Item(16, 3).Value = 140
$ws.Cells.Item(16, 4).Value = 479
$ws.Cells.Item(16, 5).Value = 7846
$ws.Cells.Item(16, 6).Value = 187
$ws.Cells.Item(16, 7).Value = 18
$ws.Cells.Item(16, 8).Value = 86

# Row 32: Rumania
$ws.Cells.Item(32, 1).Value = "Rumania"
$ws.Cells.Item(32, 2).Value = 1760
$ws.Cells.Item(32, 3).Value = 308
$ws.Cells.Item(32, 4).Value = 169
$ws.Cells.Item(32, 5).Value = 1553
$ws.Cells.Item(32, 6).Value = 34
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(32, 8).Value = 38

# Row 33: Polonia
$ws.Cells.Item(33, 1).Value = "Polonia"
$ws.Cells.Item(33, 2).Value = 1717
$ws.Cells.Item(33, 3).Value = 79
$ws.Cells.Item(33, 4).Value = 7
$ws.Cells.Item(33, 5).Value = 1691
$ws.Cells.Item(33, 6).Value = 3
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(33, 8).Value = 19

# Row 34: Japon
$ws.Cells.Item(34, 1).Value = "Japon"
$ws.Cells.Item(34, 2).Value = 1693
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 404
$ws.Cells.Item(34, 5).Value = 1237
$ws.Cells.Item(34, 6).Value = 56
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 52

# Row 35: Rusia
$ws.Cells.Item(35, 1).Value = "Rusia"
$ws.Cells.Item(35, 2).Value = 1534
$ws.Cells.Item(35, 3).Value = 270
$ws.Cells.Item(35, 4).Value = 64
$ws.Cells.Item(35, 5).Value = 1462
$ws.Cells.Item(35, 6).Value = 8
$ws.Cells.Item(35, 7).Value = 4
$ws.Cells.Item(35, 8).Value = 8

# Row 36: Pakistan
$ws.Cells.Item(36, 1).Value = "Pakistan"
$ws.Cells.Item(36, 2).Value = 1526
$ws.Cells.Item(36, 3).Value = 31
$ws.Cells.Item(36, 4).Value = 29
$ws.Cells.Item(36, 5).Value = 1484
$ws.Cells.Item(36, 6).Value = 11
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 13

# Row 50: Eslovenia
$ws.Cells.Item(50, 1).Value = "Eslovenia"
$ws.Cells.Item(50, 2).Value = 730
$ws.Cells.Item(50, 3).Value = 46
$ws.Cells.Item(50, 4).Value = 10
$ws.Cells.Item(50, 5).Value = 709
$ws.Cells.Item(50, 6).Value = 23
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 8).Value = 11

# Row 51: Republica Dominicana
$ws.Cells.Item(51, 1).Value = "Republica Dominicana"
$ws.Cells.Item(51, 2).Value = 719
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 3
$ws.Cells.Item(51, 5).Value = 688
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 28

# Row 52: Crucero
$ws.Cells.Item(52, 1).Value = "Crucero"
$ws.Cells.Item(52, 2).Value = 712
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 597
$ws.Cells.Item(52, 5).Value = 105
$ws.Cells.Item(52, 6).Value = 15
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 10

# Row 61: Irak
$ws.Cells.Item(61, 1).Value = "Irak"
$ws.Cells.Item(61, 2).Value = 547
$ws.Cells.Item(61, 3).Value = 41
$ws.Cells.Item(61, 4).Value = 143
$ws.Cells.Item(61, 5).Value = 362
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 42

# Row 62: Nueva Zelanda
$ws.Cells.Item(62, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(62, 2).Value = 514
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 56
$ws.Cells.Item(62, 5).Value = 457
$ws.Cells.Item(62, 6).Value = 1
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 1

# Row 66: Libano
$ws.Cells.Item(66, 1).Value = "Libano"
$ws.Cells.Item(66, 2).Value = 438
$ws.Cells.Item(66, 3).Value = 26
$ws.Cells.Item(66, 4).Value = 30
$ws.Cells.Item(66, 5).Value = 398
$ws.Cells.Item(66, 6).Value = 4
$ws.Cells.Item(66, 7).Value = 2
$ws.Cells.Item(66, 8).Value = 10

# Row 67: Lituania
$ws.Cells.Item(67, 1).Value = "Lituania"
$ws.Cells.Item(67, 2).Value = 437
$ws.Cells.Item(67, 3).Value = 43
$ws.Cells.Item(67, 4).Value = 1
$ws.Cells.Item(67, 5).Value = 429
$ws.Cells.Item(67, 6).Value = 2
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 7

# Row 68: Marruecos
$ws.Cells.Item(68, 1).Value = "Marruecos"
$ws.Cells.Item(68, 2).Value = 437
$ws.Cells.Item(68, 3).Value = 35
$ws.Cells.Item(68, 4).Value = 12
$ws.Cells.Item(68, 5).Value = 399
$ws.Cells.Item(68, 6).Value = 1
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 26

# Row 69: Armenia
$ws.Cells.Item(69, 1).Value = "Armenia"
$ws.Cells.Item(69, 2).Value = 424
$ws.Cells.Item(69, 3).Value = 17
$ws.Cells.Item(69, 4).Value = 30
$ws.Cells.Item(69, 5).Value = 391
$ws.Cells.Item(69, 6).Value = 6
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 3

# Row 70: Ucrania
$ws.Cells.Item(70, 1).Value = "Ucrania"
$ws.Cells.Item(70, 2).Value = 418
$ws.Cells.Item(70, 3).Value = 62
$ws.Cells.Item(70, 4).Value = 5
$ws.Cells.Item(70, 5).Value = 404
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 9

# Row 75: Eslovaquia
$ws.Cells.Item(75, 1).Value = "Eslovaquia"
$ws.Cells.Item(75, 2).Value = 314
$ws.Cells.Item(75, 3).Value = 22
$ws.Cells.Item(75, 4).Value = 2
$ws.Cells.Item(75, 5).Value = 312
$ws.Cells.Item(75, 6).Value = 1
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0

# Row 76: Principado de Andorra
$ws.Cells.Item(76, 1).Value = "Principado de Andorra"
$ws.Cells.Item(76, 2).Value = 308
$ws.Cells.Item(76, 3).Value = 0
$ws.Cells.Item(76, 4).Value = 1
$ws.Cells.Item(76, 5).Value = 304
$ws.Cells.Item(76, 6).Value = 10
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 3

# Row 77: Uruguay
$ws.Cells.Item(77, 1).Value = "Uruguay"
$ws.Cells.Item(77, 2).Value = 304
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 303
$ws.Cells.Item(77, 6).Value = 9
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 1

# Row 78: Taiwan
$ws.Cells.Item(78, 1).Value = "Taiwan"
$ws.Cells.Item(78, 2).Value = 298
$ws.Cells.Item(78, 3).Value = 15
$ws.Cells.Item(78, 4).Value = 30
$ws.Cells.Item(78, 5).Value = 266
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 2

# Row 79: Costa Rica
$ws.Cells.Item(79, 1).Value = "Costa Rica"
$ws.Cells.Item(79, 2).Value = 295
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(79, 4).Value = 3
$ws.Cells.Item(79, 5).Value = 290
$ws.Cells.Item(79, 6).Value = 6
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 2

# Row 81: Kuwait
$ws.Cells.Item(81, 1).Value = "Kuwait"
$ws.Cells.Item(81, 2).Value = 255
$ws.Cells.Item(81, 3).Value = 20
$ws.Cells.Item(81, 4).Value = 67
$ws.Cells.Item(81, 5).Value = 188
$ws.Cells.Item(81, 6).Value = 12
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 0

# Row 82: Kazajistan
$ws.Cells.Item(82, 1).Value = "Kazajistan"
$ws.Cells.Item(82, 2).Value = 251
$ws.Cells.Item(82, 3).Value = 23
$ws.Cells.Item(82, 4).Value = 18
$ws.Cells.Item(82, 5).Value = 232
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 1

# Row 83: Jordania
$ws.Cells.Item(83, 1).Value = "Jordania"
$ws.Cells.Item(83, 2).Value = 246
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 18
$ws.Cells.Item(83, 5).Value = 227
$ws.Cells.Item(83, 6).Value = 3
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 1

# Row 84: Republica de Macedonia
$ws.Cells.Item(84, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(84, 2).Value = 241
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = 3
$ws.Cells.Item(84, 5).Value = 234
$ws.Cells.Item(84, 6).Value = 1
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 4

# Row 87: Albania
$ws.Cells.Item(87, 1).Value = "Albania"
$ws.Cells.Item(87, 2).Value = 212
$ws.Cells.Item(87, 3).Value = 15
$ws.Cells.Item(87, 4).Value = 33
$ws.Cells.Item(87, 5).Value = 169
$ws.Cells.Item(87, 6).Value = 3
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 10

# Row 88: Burkina Faso
$ws.Cells.Item(88, 1).Value = "Burkina Faso"
$ws.Cells.Item(88, 2).Value = 207
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 21
$ws.Cells.Item(88, 5).Value = 175
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 11

# Row 97: Costa de Marfil
$ws.Cells.Item(97, 1).Value = "Costa de Marfil"
$ws.Cells.Item(97, 2).Value = 140
$ws.Cells.Item(97, 3).Value = 39
$ws.Cells.Item(97, 4).Value = 3
$ws.Cells.Item(97, 5).Value = 137
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0

# Row 98: Uzbekistan
$ws.Cells.Item(98, 1).Value = "Uzbekistan"
$ws.Cells.Item(98, 2).Value = 133
$ws.Cells.Item(98, 3).Value = 29
$ws.Cells.Item(98, 4).Value = 7
$ws.Cells.Item(98, 5).Value = 124
$ws.Cells.Item(98, 6).Value = 8
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 2

# Row 99: Senegal
$ws.Cells.Item(99, 1).Value = "Senegal"
$ws.Cells.Item(99, 2).Value = 130
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 18
$ws.Cells.Item(99, 5).Value = 112
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0

# Row 100: Brunei
$ws.Cells.Item(100, 1).Value = "Brunei"
$ws.Cells.Item(100, 2).Value = 126
$ws.Cells.Item(100, 3).Value = 6
$ws.Cells.Item(100, 4).Value = 34
$ws.Cells.Item(100, 5).Value = 91
$ws.Cells.Item(100, 6).Value = 1
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 1

# Row 101: Cuba
$ws.Cells.Item(101, 1).Value = "Cuba"
$ws.Cells.Item(101, 2).Value = 119
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 4
$ws.Cells.Item(101, 5).Value = 112
$ws.Cells.Item(101, 6).Value = 2
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 3

# Row 102: Venezuela
$ws.Cells.Item(102, 1).Value = "Venezuela"
$ws.Cells.Item(102, 2).Value = 119
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(102, 4).Value = 39
$ws.Cells.Item(102, 5).Value = 78
$ws.Cells.Item(102, 6).Value = 2
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 2

# Row 103: Sri Lanka
$ws.Cells.Item(103, 1).Value = "Sri Lanka"
$ws.Cells.Item(103, 2).Value = 115
$ws.Cells.Item(103, 3).Value = 2
$ws.Cells.Item(103, 4).Value = 10
$ws.Cells.Item(103, 5).Value = 104
$ws.Cells.Item(103, 6).Value = 5
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 1

# Row 104: Honduras
$ws.Cells.Item(104, 1).Value = "Honduras"
$ws.Cells.Item(104, 2).Value = 110
$ws.Cells.Item(104, 3).Value = 15
$ws.Cells.Item(104, 4).Value = 3
$ws.Cells.Item(104, 5).Value = 106
$ws.Cells.Item(104, 6).Value = 4
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 1

# Row 105: Afganistan
$ws.Cells.Item(105, 1).Value = "Afganistan"
$ws.Cells.Item(105, 2).Value = 110
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 2
$ws.Cells.Item(105, 5).Value = 104
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 4

# Row 106: Estado de Palestina
$ws.Cells.Item(106, 1).Value = "Estado de Palestina"
$ws.Cells.Item(106, 2).Value = 106
$ws.Cells.Item(106, 3).Value = 2
$ws.Cells.Item(106, 4).Value = 18
$ws.Cells.Item(106, 5).Value = 87
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 1

# Row 107: Camboya
$ws.Cells.Item(107, 1).Value = "Camboya"
$ws.Cells.Item(107, 2).Value = 103
$ws.Cells.Item(107, 3).Value = 4
$ws.Cells.Item(107, 4).Value = 21
$ws.Cells.Item(107, 5).Value = 82
$ws.Cells.Item(107, 6).Value = 1
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 0

# Row 108: Mauricio
$ws.Cells.Item(108, 1).Value = "Mauricio"
$ws.Cells.Item(108, 2).Value = 102
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 0
$ws.Cells.Item(108, 5).Value = 100
$ws.Cells.Item(108, 6).Value = 1
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 2

# Row 109: Guadalupe
$ws.Cells.Item(109, 1).Value = "Guadalupe"
$ws.Cells.Item(109, 2).Value = 102
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 17
$ws.Cells.Item(109, 5).Value = 83
$ws.Cells.Item(109, 6).Value = 4
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 2

# Row 131: Isla de Man
$ws.Cells.Item(131, 1).Value = "Isla de Man"
$ws.Cells.Item(131, 2).Value = 37
$ws.Cells.Item(131, 3).Value = 5
$ws.Cells.Item(131, 4).Value = 0
$ws.Cells.Item(131, 5).Value = 37
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 0

# Row 132: Macao
$ws.Cells.Item(132, 1).Value = "Macao"
$ws.Cells.Item(132, 2).Value = 34
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 10
$ws.Cells.Item(132, 5).Value = 24
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 0

# Row 133: Guatemala
$ws.Cells.Item(133, 1).Value = "Guatemala"
$ws.Cells.Item(133, 2).Value = 34
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 10
$ws.Cells.Item(133, 5).Value = 23
$ws.Cells.Item(133, 6).Value = 1
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 1

